# feat: Add console menu
# Rearranges the cashier/counter assignment schedule on Sheet1: the staff
# name + shift-time entries in columns B:E are reshuffled across rows,
# while the row labels in column A (Caja ...) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E values per row (row label in column A is unchanged).
$data = @{
    1  = @("CHAVEZ ONOFRE, CAMILA GERALDINE - 08:15AM - 05:00PM", "HEREDIA CAHUAYA, SUSAN NAYELLI - 05:45PM - 10:45PM", "", "")
    2  = @("MEZA MELO, NORMA FERNANDA - 08:30AM - 12:15PM", "ZAVALA SOSA, NICOLE - 01:00PM - 10:00PM", "", "")
    3  = @("FLORES PAREDES, LOURDES - 07:00AM - 09:45AM", "Del Aguila Murayari, Darla - 09:45AM - 01:30PM", "VARGAS CASTRO, LOANA VICTORIA - 02:45PM - 06:30PM", "LA ROSA EUSEBIO, SHADIA SHAMIRA - 07:00PM - 10:45PM")
    4  = @("BRICEÑO LUNA, JESSICA ARACELI - 07:45AM - 06:45PM", "BONILLA SÁNCHEZ, RAÚL FERNANDO - 07:00PM - 10:45PM", "", "")
    5  = @("ALVITE CORNEJO, ANGIE LUCERO - 08:30AM - 12:15PM", "TITO LAURA, NANCY FIORELLA - 12:30PM - 09:30PM", "", "")
    6  = @("HUAMANI TORRES, LUIS RODRIGO - 06:30AM - 03:15PM", "AYQUIPA MONTENEGRO, VALERIA ESTEFANY - 05:00PM - 08:45PM", "", "")
    7  = @("ERIQUE CALLE, MARIA ANTONIETA - 10:30AM - 07:15PM", "", "", "")
    8  = @("BARRIENTOS JERI, MILAGROS NICOL - 03:00PM - 06:45PM", "CHIARA LIMA, AUGUSTO SEBASTIAN - 07:00PM - 10:45PM", "", "")
    9  = @("DUEÑAS QUISPE, JUDYTH EVELYN - 09:00AM - 12:45PM", "ZEVALLOS ZANCA, VERONICA LUZ - 02:00PM - 11:00PM", "", "")
    10 = @("QUIQUIA MALLQUI, CYNTHIA ANGELLINE - 07:30AM - 11:15AM", "SOTO VELAZCO, EMIR ALESSANDRO - 11:15AM - 03:00PM", "RIVERA RAZA, CATHERINE - 05:00PM - 08:45PM", "")
    11 = @("YANQUI BRAVO, MIRIAN LUZ - 08:45AM - 12:30PM", "YACILA GRANDEZ, RODRIGO ANDRE - 02:00PM - 05:45PM", "NORABUENA UCHUYA, VALERIA SOFIA - 05:45PM - 09:30PM", "")
    12 = @("HUAMAN HUAMANI, ALEXIS JAVIER - 09:30AM - 01:15PM", "ALTAMIZA MATOS, MERYEIN - 02:00PM - 05:45PM", "CARDENAS RICAPA, FABRIZIO ESTEBAN - 06:00PM - 09:45PM", "")
    13 = @("SUAREZ JARA, YENNIFER YUSSARA - 09:30AM - 01:15PM", "BRENIS LÁRTIGA, SEBASTIÁN - 02:00PM - 05:45PM", "LAVADO LAZARO, CELIA ELIZABETH - 06:00PM - 09:45PM", "")
    14 = @("RUIZ SANTOS, CIELO CRISTHINA - 09:45AM - 01:30PM", "HUAYNATES ALTAMIRANO, JIM HANS - 03:45PM - 07:30PM", "", "")
    15 = @("VEGA RIVAS, ANDREA FERNANDA - 10:00AM - 01:45PM", "AYALA TAPIA, DARCIE SOL - 05:00PM - 08:45PM", "", "")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
}
